$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 previously had a trailing "Header_004_r_001_c_005" cell in column G
# that the updated test data file no longer has (the shared string that used
# to live there became unused once the new TestTable_005..009 blocks were
# appended to the sheet).
$ws.Range("G18").ClearContents()

# Append five more sample "TestTable" blocks (TestTable_005 .. TestTable_009),
# each one a table title row followed by a header row and several data rows,
# mirroring the TestTable_003 / TestTable_004 blocks already on the sheet.
$newCells = [ordered]@{
    "B23" = "TestTable_005"
    "C24" = "Header_005_r_001_c_001"
    "D24" = "Header_005_r_001_c_002"
    "E24" = "Header_005_r_001_c_003"
    "F24" = "Header_005_r_001_c_004"
    "G24" = "Header_005_r_001_c_005"
    "D25" = "Data_005_r_001_c_002"
    "E25" = "Data_005_r_001_c_003"
    "F25" = "Data_005_r_001_c_004"
    "G25" = "Data_005_r_001_c_005"
    "C26" = "Data_005_r_002_c_001"
    "D26" = "Data_005_r_002_c_002"
    "E26" = "Data_005_r_002_c_003"
    "F26" = "Data_005_r_002_c_004"
    "G26" = "Data_005_r_002_c_005"
    "C27" = "Data_005_r_003_c_001"
    "D27" = "Data_005_r_003_c_002"
    "E27" = "Data_005_r_003_c_003"
    "F27" = "Data_005_r_003_c_004"
    "G27" = "Data_005_r_003_c_005"
    "C28" = "Data_005_r_004_c_001"
    "D28" = "Data_005_r_004_c_002"
    "E28" = "Data_005_r_004_c_003"
    "F28" = "Data_005_r_004_c_004"
    "G28" = "Data_005_r_004_c_005"
    "B29" = "TestTable_006"
    "C30" = "Header_005_r_001_c_001"
    "D30" = "Header_005_r_001_c_002"
    "E30" = "Header_005_r_001_c_003"
    "F30" = "Header_005_r_001_c_004"
    "G30" = "Header_005_r_001_c_005"
    "C31" = "Data_005_r_001_c_001"
    "D31" = "Data_005_r_001_c_002"
    "E31" = "Data_005_r_001_c_003"
    "F31" = "Data_005_r_001_c_004"
    "C32" = "Data_005_r_002_c_001"
    "D32" = "Data_005_r_002_c_002"
    "E32" = "Data_005_r_002_c_003"
    "F32" = "Data_005_r_002_c_004"
    "G32" = "Data_005_r_002_c_005"
    "C33" = "Data_005_r_003_c_001"
    "D33" = "Data_005_r_003_c_002"
    "E33" = "Data_005_r_003_c_003"
    "F33" = "Data_005_r_003_c_004"
    "G33" = "Data_005_r_003_c_005"
    "C34" = "Data_005_r_004_c_001"
    "D34" = "Data_005_r_004_c_002"
    "E34" = "Data_005_r_004_c_003"
    "F34" = "Data_005_r_004_c_004"
    "G34" = "Data_005_r_004_c_005"
    "B35" = "TestTable_007"
    "C36" = "Header_005_r_001_c_001"
    "D36" = "Header_005_r_001_c_002"
    "E36" = "Header_005_r_001_c_003"
    "F36" = "Header_005_r_001_c_004"
    "G36" = "Header_005_r_001_c_005"
    "C37" = "Data_005_r_001_c_001"
    "D37" = "Data_005_r_001_c_002"
    "E37" = "Data_005_r_001_c_003"
    "F37" = "Data_005_r_001_c_004"
    "G37" = "Data_005_r_001_c_005"
    "C38" = "Data_005_r_002_c_001"
    "D38" = "Data_005_r_002_c_002"
    "E38" = "Data_005_r_002_c_003"
    "F38" = "Data_005_r_002_c_004"
    "G38" = "Data_005_r_002_c_005"
    "C39" = "Data_005_r_003_c_001"
    "D39" = "Data_005_r_003_c_002"
    "E39" = "Data_005_r_003_c_003"
    "F39" = "Data_005_r_003_c_004"
    "G39" = "Data_005_r_003_c_005"
    "D40" = "Data_005_r_004_c_002"
    "E40" = "Data_005_r_004_c_003"
    "F40" = "Data_005_r_004_c_004"
    "G40" = "Data_005_r_004_c_005"
    "B41" = "TestTable_008"
    "C42" = "Header_005_r_001_c_001"
    "D42" = "Header_005_r_001_c_002"
    "E42" = "Header_005_r_001_c_003"
    "F42" = "Header_005_r_001_c_004"
    "G42" = "Header_005_r_001_c_005"
    "C43" = "Data_005_r_001_c_001"
    "D43" = "Data_005_r_001_c_002"
    "E43" = "Data_005_r_001_c_003"
    "F43" = "Data_005_r_001_c_004"
    "G43" = "Data_005_r_001_c_005"
    "C44" = "Data_005_r_002_c_001"
    "D44" = "Data_005_r_002_c_002"
    "E44" = "Data_005_r_002_c_003"
    "F44" = "Data_005_r_002_c_004"
    "G44" = "Data_005_r_002_c_005"
    "C45" = "Data_005_r_003_c_001"
    "D45" = "Data_005_r_003_c_002"
    "E45" = "Data_005_r_003_c_003"
    "F45" = "Data_005_r_003_c_004"
    "G45" = "Data_005_r_003_c_005"
    "C46" = "Data_005_r_004_c_001"
    "D46" = "Data_005_r_004_c_002"
    "E46" = "Data_005_r_004_c_003"
    "F46" = "Data_005_r_004_c_004"
    "B47" = "TestTable_009"
    "C48" = "Header_005_r_001_c_001"
    "D48" = "Header_005_r_001_c_002"
    "E48" = "Header_005_r_001_c_003"
    "F48" = "Header_005_r_001_c_004"
    "G48" = "Header_005_r_001_c_005"
    "C49" = "Data_005_r_001_c_001"
    "D49" = "Data_005_r_001_c_002"
    "E49" = "Data_005_r_001_c_003"
    "F49" = "Data_005_r_001_c_004"
    "G49" = "Data_005_r_001_c_005"
    "C50" = "Data_005_r_002_c_001"
    "D50" = "Data_005_r_002_c_002"
    "F50" = "Data_005_r_002_c_004"
    "G50" = "Data_005_r_002_c_005"
    "C51" = "Data_005_r_003_c_001"
    "D51" = "Data_005_r_003_c_002"
    "E51" = "Data_005_r_003_c_003"
    "F51" = "Data_005_r_003_c_004"
    "G51" = "Data_005_r_003_c_005"
    "C52" = "Data_005_r_004_c_001"
    "D52" = "Data_005_r_004_c_002"
    "E52" = "Data_005_r_004_c_003"
    "F52" = "Data_005_r_004_c_004"
    "G52" = "Data_005_r_004_c_005"
}

foreach ($addr in $newCells.Keys) {
    $ws.Range($addr).Value = $newCells[$addr]
}

# Match the workbook's recorded selection after the edit.
$ws.Range("G19").Select() | Out-Null
